$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 297.76923
$ws.Range("I53").Value = 303.16666
$ws.Range("J53").Value = 293.14285
$ws.Range("K53").Value = 303.16666
$ws.Range("L53").Value = 293.14285
$ws.Range("M53").Value = 333.83334
$ws.Range("N53").Value = -1567.14285
$ws.Range("H62").Value = 4549.75
$ws.Range("I62").Value = 3399.6667
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3399.6667
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -2775.6667
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 4549.75
$ws.Range("I65").Value = 3399.6667
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 16998.3335
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -13878.3335
$ws.Range("N65").Value = -46240
$ws.Range("H88").Value = 4655.3335
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 1983
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 1983
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -2795
$ws.Range("H91").Value = 4655.3335
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 1983
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 1983
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -4791
$ws.Range("H111").Value = 2164.875
$ws.Range("I111").Value = 964.5
$ws.Range("J111").Value = 5766
$ws.Range("K111").Value = 2893.5
$ws.Range("L111").Value = 17298
$ws.Range("M111").Value = 173.5
$ws.Range("N111").Value = -23432
$ws.Range("H131").Value = 2715.8333
$ws.Range("I131").Value = 1098.3334
$ws.Range("J131").Value = 4333.3335
$ws.Range("K131").Value = 3295.0002
$ws.Range("L131").Value = 13000.0005
$ws.Range("M131").Value = 1744.9998
$ws.Range("N131").Value = -23080.0005
$ws.Range("H137").Value = 3113.889
$ws.Range("I137").Value = 2914.7058
$ws.Range("J137").Value = 6500
$ws.Range("K137").Value = 8744.117400000001
$ws.Range("L137").Value = 19500
$ws.Range("M137").Value = -6194.117400000001
$ws.Range("N137").Value = -24600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 26997.5
$ws.Range("I31").Value = 3995
$ws.Range("J31").Value = 50000
$ws.Range("K31").Value = 3995
$ws.Range("L31").Value = 50000
$ws.Range("M31").Value = -3701
$ws.Range("N31").Value = -50588
$ws.Range("H32").Value = 7457.3447
$ws.Range("I32").Value = 7457.3447
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7457.3447
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7170.3447
$ws.Range("H63").Value = 7002
$ws.Range("I63").Value = 500
$ws.Range("J63").Value = 20006
$ws.Range("K63").Value = 500
$ws.Range("L63").Value = 20006
$ws.Range("M63").Value = 186
$ws.Range("N63").Value = -21378
$ws.Range("H66").Value = 7002
$ws.Range("I66").Value = 500
$ws.Range("J66").Value = 20006
$ws.Range("K66").Value = 2500
$ws.Range("L66").Value = 100030
$ws.Range("M66").Value = 932
$ws.Range("N66").Value = -106894
$ws.Range("H132").Value = 2032.24
$ws.Range("I132").Value = 1490.9048
$ws.Range("J132").Value = 4874.25
$ws.Range("K132").Value = 4472.7144
$ws.Range("L132").Value = 14622.75
$ws.Range("M132").Value = -1942.7144
$ws.Range("N132").Value = -19682.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7684.5
$ws.Range("I20").Value = 1527
$ws.Range("J20").Value = 19999.5
$ws.Range("K20").Value = 1527
$ws.Range("L20").Value = 19999.5
$ws.Range("M20").Value = -1280
$ws.Range("N20").Value = -20493.5
$ws.Range("H80").Value = 678.53845
$ws.Range("I80").Value = 205.75
$ws.Range("J80").Value = 888.6667
$ws.Range("K80").Value = 205.75
$ws.Range("L80").Value = 888.6667
$ws.Range("M80").Value = 792.25
$ws.Range("N80").Value = -2884.6667
$ws.Range("H83").Value = 678.53845
$ws.Range("I83").Value = 205.75
$ws.Range("J83").Value = 888.6667
$ws.Range("K83").Value = 1028.75
$ws.Range("L83").Value = 4443.3335
$ws.Range("M83").Value = 3963.25
$ws.Range("N83").Value = -14427.3335
$ws.Range("H86").Value = 5626.7144
$ws.Range("I86").Value = 2996.3333
$ws.Range("J86").Value = 7599.5
$ws.Range("K86").Value = 2996.3333
$ws.Range("L86").Value = 7599.5
$ws.Range("M86").Value = -1873.3333
$ws.Range("N86").Value = -9845.5
$ws.Range("H89").Value = 5626.7144
$ws.Range("I89").Value = 2996.3333
$ws.Range("J89").Value = 7599.5
$ws.Range("K89").Value = 14981.6665
$ws.Range("L89").Value = 37997.5
$ws.Range("M89").Value = -9365.666499999999
$ws.Range("N89").Value = -49229.5
$ws.Range("H107").Value = 1118.3572
$ws.Range("I107").Value = 1138.0834
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1138.0834
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 781.9166
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1365.9
$ws.Range("I132").Value = 962.2222
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2886.6666
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -356.6666
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 2850.75
$ws.Range("I134").Value = 2850.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8552.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6017.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 164.33333
$ws.Range("I2").Value = 111
$ws.Range("J2").Value = 271
$ws.Range("K2").Value = 666
$ws.Range("L2").Value = 1626
$ws.Range("M2").Value = -553
$ws.Range("N2").Value = -1852
$ws.Range("H5").Value = 1344.25
$ws.Range("I5").Value = 1344.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4032.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3920.75
$ws.Range("H40").Value = 363.2
$ws.Range("I40").Value = 333
$ws.Range("J40").Value = 383.33334
$ws.Range("K40").Value = 1332
$ws.Range("L40").Value = 1533.33336
$ws.Range("M40").Value = -1263
$ws.Range("N40").Value = -1671.33336
$ws.Range("H132").Value = 1494
$ws.Range("I132").Value = 1494
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13446
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10916
$ws.Range("H135").Value = 1344.25
$ws.Range("I135").Value = 1344.25
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12098.25
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9563.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 19998.5
$ws.Range("I22").Value = 19999
$ws.Range("J22").Value = 19998
$ws.Range("K22").Value = 19999
$ws.Range("L22").Value = 19998
$ws.Range("M22").Value = -19470
$ws.Range("N22").Value = -21056
$ws.Range("H70").Value = 1849.25
$ws.Range("I70").Value = 1799
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 1799
$ws.Range("L70").Value = 2000
$ws.Range("M70").Value = -1529
$ws.Range("N70").Value = -2540
$ws.Range("H73").Value = 1849.25
$ws.Range("I73").Value = 1799
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 1799
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = -863
$ws.Range("N73").Value = -3872
$ws.Range("H102").Value = 2473.7778
$ws.Range("I102").Value = 2783.1428
$ws.Range("J102").Value = 1391
$ws.Range("K102").Value = 2783.1428
$ws.Range("L102").Value = 1391
$ws.Range("M102").Value = -1161.1428
$ws.Range("N102").Value = -4635
$ws.Range("H126").Value = 5933.3335
$ws.Range("I126").Value = 5933.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 17800.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -15330.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2470
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -10376
$ws.Range("H55").Value = 1235.9048
$ws.Range("I55").Value = 2387
$ws.Range("J55").Value = 372.58334
$ws.Range("K55").Value = 2387
$ws.Range("L55").Value = 372.58334
$ws.Range("M55").Value = -2214
$ws.Range("N55").Value = -718.58334
$ws.Range("H132").Value = 2559.5483
$ws.Range("I132").Value = 2062.4583
$ws.Range("J132").Value = 4263.857
$ws.Range("K132").Value = 6187.374899999999
$ws.Range("L132").Value = 12791.571
$ws.Range("M132").Value = -3657.374899999999
$ws.Range("N132").Value = -17851.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 4888
$ws.Range("I70").Value = 4888
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4888
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4573
$ws.Range("H73").Value = 4888
$ws.Range("I73").Value = 4888
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4888
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3796
$ws.Range("H132").Value = 1689.091
$ws.Range("I132").Value = 953.55554
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2860.66662
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -330.66662
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 2322.2593
$ws.Range("I136").Value = 2133.5833
$ws.Range("J136").Value = 3831.6667
$ws.Range("K136").Value = 6400.749899999999
$ws.Range("L136").Value = 11495.0001
$ws.Range("M136").Value = -3850.749899999999
$ws.Range("N136").Value = -16595.0001
